$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'92.560.56"
$ws.Range("E2").Value = "  -5.77%  "
$ws.Range("D3").Value = "'3.359.04"
$ws.Range("E3").Value = "  -1.34%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'231.57"
$ws.Range("E5").Value = "  -8.94%  "
$ws.Range("D6").Value = "'628.84"
$ws.Range("E6").Value = "  -5.35%  "
$ws.Range("D7").Value = "'1.39"
$ws.Range("E7").Value = "  -8.42%  "
$ws.Range("D8").Value = "'0.389"
$ws.Range("E8").Value = "  -10.08%  "
$ws.Range("E9").Value = "  +0.07%  "
$ws.Range("D10").Value = "'0.939"
$ws.Range("E10").Value = "  -10.99%  "
$ws.Range("D11").Value = "'3.354.73"
$ws.Range("E11").Value = "  -1.38%  "
$ws.Range("D12").Value = "'0.195"
$ws.Range("E12").Value = "  -7.03%  "
$ws.Range("D13").Value = "'40.19"
$ws.Range("E13").Value = "  -11.71%  "
$ws.Range("D14").Value = "'5.98"
$ws.Range("E14").Value = "  -2.76%  "
$ws.Range("D15").Value = "'92.757.49"
$ws.Range("E15").Value = "  -5.37%  "
$ws.Range("D16").Value = "'3.985.35"
$ws.Range("E16").Value = "  -1.64%  "
$ws.Range("D17").Value = "'0.0000243"
$ws.Range("E17").Value = "  -6.11%  "
$ws.Range("D18").Value = "'7.96"
$ws.Range("E18").Value = "  -12.99%  "
$ws.Range("D19").Value = "'3.354.35"
$ws.Range("E19").Value = "  -1.40%  "
$ws.Range("D20").Value = "'16.85"
$ws.Range("E20").Value = "  -8.87%  "
$ws.Range("D21").Value = "'11.03"
$ws.Range("E21").Value = "  -2.97%  "
$ws.Range("D22").Value = "'484.52"
$ws.Range("E22").Value = "  -5.31%  "
$ws.Range("D23").Value = "'0.449"
$ws.Range("E23").Value = "  -16.95%  "
$ws.Range("D24").Value = "'3.12"
$ws.Range("E24").Value = "  -8.77%  "
$ws.Range("D25").Value = "'0.0000185"
$ws.Range("E25").Value = "  -8.29%  "
$ws.Range("D26").Value = "'6.26"
$ws.Range("E26").Value = "  -7.97%  "
$ws.Range("D27").Value = "'89.55"
$ws.Range("E27").Value = "  -8.32%  "
$ws.Range("B28").Value = "WrappedeETH"
$ws.Range("C28").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D28").Value = "'3.540.41"
$ws.Range("E28").Value = "  -1.34%  "
$ws.Range("B29").Value = "Aptos"
$ws.Range("C29").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D29").Value = "'11.47"
$ws.Range("E29").Value = "  -7.86%  "
$ws.Range("E30").Value = "  +0.11%  "
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").Value = "'11.23"
$ws.Range("E31").Value = "  -7.75%  "
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").Value = "'2.63"
$ws.Range("E32").Value = "  -8.87%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "'0.131"
$ws.Range("E33").Value = "  -8.91%  "
$ws.Range("B34").Value = "Binance-PegBSC-USD"
$ws.Range("C34").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D34").Value = "'0.999"
$ws.Range("E34").Value = "  +0.10%  "
$ws.Range("B35").Value = "Cronos"
$ws.Range("C35").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D35").Value = "'0.170"
$ws.Range("E35").Value = "  -9.56%  "
$ws.Range("B36").Value = "EthereumClassic"
$ws.Range("C36").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D36").Value = "'28.39"
$ws.Range("E36").Value = "  -2.56%  "
$ws.Range("B37").Value = "PolygonEcosystemToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D37").Value = "'0.527"
$ws.Range("E37").Value = "  -7.05%  "
$ws.Range("B38").Value = "Bittensor"
$ws.Range("C38").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D38").Value = "'532.13"
$ws.Range("E38").Value = "  +1.15%  "
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D39").Value = "'7.39"
$ws.Range("E39").Value = "  -6.90%  "
$ws.Range("E40").Value = "  +0.01%  "
$ws.Range("B41").Value = "Fetch.AI"
$ws.Range("C41").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D41").Value = "'1.39"
$ws.Range("E41").Value = "  -6.53%  "
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").Value = "'0.147"
$ws.Range("E42").Value = "  -4.60%  "
$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D43").Value = "'0.874"
$ws.Range("E43").Value = "  +1.41%  "
$ws.Range("B44").Value = "WhiteBITCoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D44").Value = "'24.01"
$ws.Range("E44").Value = "  -1.66%  "
$ws.Range("B45").Value = "MantraDAO"
$ws.Range("C45").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D45").Value = "'3.61"
$ws.Range("E45").Value = "  -2.08%  "
$ws.Range("B46").Value = "ImmutableX"
$ws.Range("C46").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D46").Value = "'1.65"
$ws.Range("E46").Value = "  -5.00%  "
$ws.Range("B47").Value = "Filecoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D47").Value = "'5.46"
$ws.Range("E47").Value = "  -3.49%  "
$ws.Range("B48").Value = "OKB"
$ws.Range("C48").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D48").Value = "'53.03"
$ws.Range("E48").Value = "  -5.29%  "
$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").Value = "'0.0392"
$ws.Range("E49").Value = "  -8.16%  "
$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").Value = "'2.12"
$ws.Range("E50").Value = "  -5.10%  "
$ws.Range("B51").Value = "dogwifhat"
$ws.Range("C51").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D51").Value = "'3.13"
$ws.Range("E51").Value = "  -3.39%  "
